# The deck's design theme (ppt/theme/theme1.xml, wired to the slide master)
# is switched from the custom "Integral / Red Violet" palette to the
# stock PowerPoint "Office Theme" palette. Re-point every slot of the
# theme's color scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink -
# ThemeColorScheme.Colors(1..12)) at the standard Office Theme RGB values.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

function HexRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 0; $i -lt $officeColors.Length; $i++) {
    $colorScheme.Colors($i + 1).RGB = HexRgb($officeColors[$i])
}
